# Refresh the cryptos price list (columns D "Price" and E "Volume(1h)")
# with updated scrape values, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.543.88"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "1.846.33"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  -1.19%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4642"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "1.855.84"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.957"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.131"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06679"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001035"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "27.542.76"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.387"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9763"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09400"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.593"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.304"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.339"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06048"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02226"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.294"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5893"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5586"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.911"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06701"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.006"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.923"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -13.80%  "
